$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 14 de Abril de 2020 a las 19:22'
$ws.Range('B6').Value = 7597
$ws.Range('C6').Value = 1240
$ws.Range('D6').Value = 5973
$ws.Range('E6').Value = 384
$ws.Range('B9').Value = 4752
$ws.Range('C9').Value = 1587
$ws.Range('D9').Value = 2718
$ws.Range('E9').Value = 447
$ws.Range('B13').Value = 3323
$ws.Range('C13').Value = 1198
$ws.Range('D13').Value = 1778
$ws.Range('E13').Value = 347
$ws.Range('B18').Value = 2672
$ws.Range('C18').Value = 896
$ws.Range('D18').Value = 1560
$ws.Range('E18').Value = 216
$ws.Range('B19').Value = 2269
$ws.Range('C19').Value = 658
$ws.Range('D19').Value = 1350
$ws.Range('E19').Value = 261
$ws.Range('A24').Value = 'Segovia'
$ws.Range('B24').Value = 1994
$ws.Range('C24').Value = 563
$ws.Range('D24').Value = 1276
$ws.Range('E24').Value = 155
$ws.Range('A25').Value = 'A Coruña'
$ws.Range('B25').Value = 1969
$ws.Range('C25').Value = 333
$ws.Range('D25').Value = 1788
$ws.Range('E25').Value = 67
$ws.Range('A26').Value = 'Leon'
$ws.Range('B26').Value = 1865
$ws.Range('C26').Value = 854
$ws.Range('D26').Value = 747
$ws.Range('E26').Value = 264
$ws.Range('A27').Value = 'Granada'
$ws.Range('B27').Value = 1864
$ws.Range('C27').Value = 375
$ws.Range('D27').Value = 1309
$ws.Range('E27').Value = 180
$ws.Range('A28').Value = 'Cantabria'
$ws.Range('B28').Value = 1796
$ws.Range('C28').Value = 323
$ws.Range('D28').Value = 1353
$ws.Range('E28').Value = 120
$ws.Range('A29').Value = 'Caceres'
$ws.Range('B29').Value = 1776
$ws.Range('C29').Value = 237
$ws.Range('D29').Value = 1276
$ws.Range('E29').Value = 263
$ws.Range('B32').Value = 1257
$ws.Range('C32').Value = 548
$ws.Range('D32').Value = 561
$ws.Range('E32').Value = 148
$ws.Range('B33').Value = 1237
$ws.Range('C33').Value = 331
$ws.Range('D33').Value = 836
$ws.Range('E33').Value = 70
$ws.Range('B37').Value = 1133
$ws.Range('C37').Value = 274
$ws.Range('D37').Value = 736
$ws.Range('E37').Value = 123
$ws.Range('B38').Value = 1053
$ws.Range('C38').Value = 257
$ws.Range('D38').Value = 709
$ws.Range('E38').Value = 87
$ws.Range('B40').Value = 966
$ws.Range('C40').Value = 378
$ws.Range('D40').Value = 488
$ws.Range('E40').Value = 100
$ws.Range('B45').Value = 645
$ws.Range('C45').Value = 192
$ws.Range('D45').Value = 401
$ws.Range('E45').Value = 52
$ws.Range('A49').Value = 'Zamora'
$ws.Range('B49').Value = 459
$ws.Range('C49').Value = 175
$ws.Range('D49').Value = 230
$ws.Range('E49').Value = 54
$ws.Range('A50').Value = 'Gran Canaria'
$ws.Range('B50').Value = 456
$ws.Range('C50').Value = 194
$ws.Range('D50').Value = 235
$ws.Range('E50').Value = 27
$ws.Range('C56').Value = 17
$ws.Range('D56').Value = 63
$ws.Range('B57').Value = 78
$ws.Range('D57').Value = 56
$ws.Range('C59').Value = 25
$ws.Range('D59').Value = 17
